{"js": "// The template had the placeholder \"{{texto_comprobacion_no_cumple}}\" split\n// across three separate runs (\"{{texto_comprobacion_\", \"no_\", \"cumple}}\").\n// Find that paragraph and rewrite its range so the placeholder text is a\n// single contiguous run again (same visual text, same run formatting).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = \"{{texto_comprobacion_no_cumple}}\";\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(target) !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // Replacing the range's text collapses the multiple runs that made up\n  // the placeholder into one run, matching the merged-run edit.\n  const range = targetParagraph.getRange();\n  range.insertText(target, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The template had the placeholder \"{{texto_comprobacion_no_cumple}}\" split\n# across three separate runs (\"{{texto_comprobacion_\", \"no_\", \"cumple}}\").\n# Use Find/Replace on the whole document story to rewrite that text back\n# into a single contiguous run (same visible text, same run formatting).\n$d = $word.ActiveDocument\n\n$target = \"{{texto_comprobacion_no_cumple}}\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $target\n$find.Replacement.Text = $target\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2)\n"}
